$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("U2").Value = 0.59
$ws.Range("V2").Value = 0.36

# Row 3 updates
$ws.Range("U3").Value = 0.21
$ws.Range("V3").Value = 0.42
$ws.Range("W3").Value = 0.19
$ws.Range("X3").Value = 0.1
$ws.Range("Y3").Value = 0.04
$ws.Range("Z3").Value = 0.02
$ws.Range("AA3").Value = 0.01

# Row 4 updates
$ws.Range("U4").Value = 0.91
$ws.Range("V4").Value = 0.05
$ws.Range("W4").Value = 0
